$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts existing rows 2..81 down to 3..82)
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# Populate the new row 2 with the new publication's data
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = "B. R. Evans, A. Lowe, A. Crawford, A. Fleming, J. S. Hosking"
$ws.Range("C2").Value = "10.5194/egusphere-2025-2886"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "EGUsphere"
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "1--28"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = "Icebergs, jigsaw puzzles and genealogy: Automated multi-generational iceberg tracking and lineage reconstruction"
$ws.Range("P2").Value = "preprint"
$ws.Range("Q2").Value = "https://egusphere.copernicus.org/preprints/2025/egusphere-2025-2886/"
$ws.Range("R2").Value = "2025"
$ws.Range("S2").Value = "2025"
